$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.915.64'
$ws.Range('E2').Value = '  -0.10%  '

$ws.Range('D3').Value = '3.004.41'
$ws.Range('E3').Value = '  +0.49%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.998'
$ws.Range('E4').Value = '  -0.25%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.97'
$ws.Range('E5').Value = '  +3.19%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '146.45'
$ws.Range('E6').Value = '  +1.00%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  -0.01%  '

$ws.Range('D8').Value = '3.005.38'
$ws.Range('E8').Value = '  +0.60%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.515'
$ws.Range('E9').Value = '  -1.54%  '

$ws.Range('E10').Value = '  +1.32%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.10'
$ws.Range('E11').Value = '  +5.36%  '

$ws.Range('E12').Value = '  +0.38%  '

$ws.Range('E13').Value = '  +0.98%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.38'
$ws.Range('E14').Value = '  +0.35%  '

$ws.Range('E15').Value = '  +3.48%  '

$ws.Range('D16').Value = '3.492.41'
$ws.Range('E16').Value = '  +0.26%  '

$ws.Range('E17').Value = '  -1.60%  '

$ws.Range('D18').Value = '61.773.92'
$ws.Range('E18').Value = '  -0.30%  '

$ws.Range('D19').Value = '2.994.90'
$ws.Range('E19').Value = '  -0.12%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '448.77'
$ws.Range('E20').Value = '  -2.15%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.05'
$ws.Range('E21').Value = '  +1.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.686'
$ws.Range('E22').Value = '  +0.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.37'
$ws.Range('E23').Value = '  -0.56%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.54'
$ws.Range('E24').Value = '  +0.63%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.04'
$ws.Range('E25').Value = '  +10.22%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.22'
$ws.Range('E26').Value = '  +0.50%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.11'
$ws.Range('E27').Value = '  -0.70%  '

$ws.Range('E28').Value = '  +0.18%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.72'
$ws.Range('E29').Value = '  +4.18%  '

$ws.Range('B30').Value = 'FirstDigitalUSD'
$ws.Range('C30').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.998'
$ws.Range('E30').Value = '  -0.43%  '

$ws.Range('B31').Value = 'NEARProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.20'
$ws.Range('E31').Value = '  +3.02%  '

$ws.Range('E32').Value = '  -0.32%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '27.28'
$ws.Range('E33').Value = '  -2.93%  '

$ws.Range('E34').Value = '  +2.72%  '

$ws.Range('D35').Value = '0.0₃0837'
$ws.Range('E35').Value = '  +5.99%  '

$ws.Range('E36').Value = '  +0.49%  '

$ws.Range('E37').Value = '  +1.12%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.52'
$ws.Range('E38').Value = '  +0.72%  '

$ws.Range('E39').Value = '  -2.30%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.97'
$ws.Range('E40').Value = '  -1.71%  '

$ws.Range('E41').Value = '  +8.22%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.89'
$ws.Range('E42').Value = '  +0.97%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '398.25'
$ws.Range('E43').Value = '  +1.75%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '40.55'
$ws.Range('E44').Value = '  +10.72%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.272'
$ws.Range('E45').Value = '  +0.06%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0351'
$ws.Range('E46').Value = '  -1.02%  '

$ws.Range('D47').Value = '2.708.13'
$ws.Range('E47').Value = '  -0.55%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '131.48'
$ws.Range('E48').Value = '  +2.73%  '

$ws.Range('E49').Value = '  +0.11%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.18'
$ws.Range('E50').Value = '  +0.10%  '

$ws.Range('E51').Value = '  -1.12%  '

Write-Output "Updated cryptos list"